$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "63.542.85"
$ws.Range("E2").Value = "  +1.61%  "

Set-TextValue "D3" "3.176.41"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "593.11"
$ws.Range("E5").Value = "  -0.17%  "

Set-TextValue "D6" "136.03"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("E7").Value = "  +0.07%  "

Set-TextValue "D8" "3.173.03"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  +1.62%  "

Set-TextValue "D10" "0.143"
$ws.Range("E10").Value = "  -1.04%  "

Set-TextValue "D11" "5.36"
$ws.Range("E11").Value = "  -0.26%  "

Set-TextValue "D12" "0.456"
$ws.Range("E12").Value = "  +0.32%  "

Set-TextValue "D13" "0.0000240"
$ws.Range("E13").Value = "  -0.28%  "

Set-TextValue "D14" "34.76"
$ws.Range("E14").Value = "  +3.32%  "

Set-TextValue "D15" "3.699.67"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("E16").Value = "  -0.35%  "

Set-TextValue "D17" "3.179.77"
$ws.Range("E17").Value = "  -0.40%  "

Set-TextValue "D18" "63.499.19"
$ws.Range("E18").Value = "  +1.32%  "

Set-TextValue "D19" "6.55"
$ws.Range("E19").Value = "  -2.44%  "

Set-TextValue "D20" "462.14"
$ws.Range("E20").Value = "  -0.46%  "

Set-TextValue "D21" "13.97"
$ws.Range("E21").Value = "  -0.38%  "

Set-TextValue "D22" "0.698"
$ws.Range("E22").Value = "  -2.21%  "

Set-TextValue "D23" "7.69"
$ws.Range("E23").Value = "  +0.02%  "

Set-TextValue "D24" "13.27"
$ws.Range("E24").Value = "  -2.10%  "

Set-TextValue "D25" "83.16"
$ws.Range("E25").Value = "  -0.45%  "

Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  -1.17%  "

Set-TextValue "D29" "6.89"
$ws.Range("E29").Value = "  -0.60%  "

Set-TextValue "D30" "7.76"
$ws.Range("E30").Value = "  -1.98%  "

Set-TextValue "D31" "2.07"
$ws.Range("E31").Value = "  -0.77%  "

Set-TextValue "D32" "27.39"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("E33").Value = "  -1.33%  "

Set-TextValue "D34" "2.42"
$ws.Range("E34").Value = "  -1.14%  "

Set-TextValue "D35" "1.02"
$ws.Range("E35").Value = "  -1.85%  "

Set-TextValue "D36" "5.90"
$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D37" "51.59"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D38" "0.0₃0732"
$ws.Range("E38").Value = "  +4.45%  "

Set-TextValue "D39" "0.0392"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("E41").Value = "  -2.06%  "

Set-TextValue "D42" "2.66"
$ws.Range("E42").Value = "  +0.66%  "

Set-TextValue "D43" "394.43"
$ws.Range("E43").Value = "  -6.09%  "

Set-TextValue "D44" "2.794.60"
$ws.Range("E44").Value = "  -7.09%  "

Set-TextValue "D45" "0.252"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D46" "35.87"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D47" "127.17"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D49" "2.12"
$ws.Range("E49").Value = "  -1.63%  "

Set-TextValue "D50" "25.31"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("E51").Value = "  -0.62%  "

